# Reorder worksheets: move "总计" to be before "2022-Q2"
$wb = $excel.ActiveWorkbook

$summarySheet = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

$summarySheet.Move($q2Sheet)

# Keep "2022-Q2" as the active/selected tab, matching original state
$excel.Sheets.Item("2022-Q2").Activate()
